$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-05 12:42:11"

# Refresh the "last seen" timestamp on all existing rows (2-7)
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# Insert a new row above row 7, pushing the current row 7 down to row 8
$ws.Rows.Item(7).Insert()

# Fill in the newly inserted row 7 with the new job listing
$ws.Cells.Item(7, 1).Value = $newTimestamp
$ws.Cells.Item(7, 2).Value = "【急募】サーバー移転のプロフェッショナルを探しています!"
$ws.Cells.Item(7, 3).Value = "システム開発"
$ws.Cells.Item(7, 4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells.Item(7, 5).Value = "期限情報なし"
$ws.Hyperlinks.Add($ws.Cells.Item(7, 6), "https://www.lancers.jp/work/detail/5407189") | Out-Null
$ws.Cells.Item(7, 7).Value = 10

$wb.Save()
